# "Modified for Private Subnet"
#
# Applies the following changes to slide 1 of the Jumpbox-Architecture deck:
#   1. Merge the "Virtual " / "private cloud" runs into a single
#      "Virtual private cloud" run (TextBox 35).
#   2. Change "Public subnet" -> "Private subnet" in the VPC textbox that
#      lives inside Group 64 (TextBox 37), and shrink the textbox to the
#      single-line autofit height that results from the shorter wording.
#   3. Merge the " ssh -W %h:%p joe@" / "50.23.28.92" / " -" runs in the
#      ProxyCommand textbox (TextBox 96, inside Group 97) into one run.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------------
# 1) "Virtual " + "private cloud" -> "Virtual private cloud"
# ---------------------------------------------------------------------------
$vpcCloud = $s.Shapes.Item("TextBox 35")
$vpcCloudRange = $vpcCloud.TextFrame.TextRange
$vpcCloudRange.Characters(1, $vpcCloudRange.Length).Text = "Virtual private cloud"

# ---------------------------------------------------------------------------
# 2) "VPC Public subnet" -> "VPC Private subnet" (+ textbox height tweak)
# ---------------------------------------------------------------------------
$subnetGroup = $s.Shapes.Item("Group 64")
$subnetBox = $subnetGroup.GroupItems.Item("TextBox 37")
$subnetRange = $subnetBox.TextFrame.TextRange
$fullSubnetText = $subnetRange.Text
$subnetWordStart = $fullSubnetText.IndexOf("Public subnet") + 1
$subnetRange.Characters($subnetWordStart, "Public subnet".Length).Text = "Private subnet"
# The textbox has spAutoFit; after retyping the now one-line content,
# PowerPoint recomputes the shape height to fit the single line of text.
$subnetBox.Height = 11.886929133858267

# ---------------------------------------------------------------------------
# 3) " ssh -W %h:%p joe@" + "50.23.28.92" + " -" -> single run
# ---------------------------------------------------------------------------
$sshGroup = $s.Shapes.Item("Group 97")
$sshBox = $sshGroup.GroupItems.Item("TextBox 96")
$sshRange = $sshBox.TextFrame.TextRange
$fullSshText = $sshRange.Text
$sshStart = $fullSshText.IndexOf(" ssh -W %h:%p joe@") + 1
$sshOldLen = (" ssh -W %h:%p joe@" + "50.23.28.92" + " -").Length
$sshRange.Characters($sshStart, $sshOldLen).Text = " ssh -W %h:%p joe@50.23.28.92 -"
